# Safety Board Automation project
#
# "Questions count and category": the per-question `count` values (column B)
# get redistributed across the rows within each category (the totals per
# category - and therefore the "Category count" sheet - stay the same).
#
# "Actions list": the action titles (column A) get reshuffled into a new
# row order; the same set of action titles is kept, just reordered.

$wb = $excel.ActiveWorkbook

# --- "Questions count and category": update the count (column B) for the
#     rows whose value changed ---
$ws1 = $wb.Worksheets.Item("Questions count and category")

$ws1.Range("B2").Value = 1
$ws1.Range("B4").Value = 4
$ws1.Range("B5").Value = 1
$ws1.Range("B6").Value = 3
$ws1.Range("B8").Value = 6
$ws1.Range("B9").Value = 1
$ws1.Range("B12").Value = 2
$ws1.Range("B16").Value = 5
$ws1.Range("B17").Value = 1

# --- "Actions list": reorder the action titles in column A ---
$ws3 = $wb.Worksheets.Item("Actions list")

$ws3.Range("A3").Value = 'Remover ASINs dos bin P-2-S128M341 e P-2-S144k646 que são maiores que o bin e fazer transferência'
$ws3.Range("A4").Value = 'Refazer demarcações dos buffers ao lado da linha 8 com tamanho para suportar carrinhos de Rebin,'
$ws3.Range("A5").Value = 'Fix electrical issues and evaluate the best option to support process'
$ws3.Range("A6").Value = 'Mencionar no Stand up a dica de segurança """" Sempre que não estiver fazendo uso do paleteira manual ela deve estar em baixo de um palete e com a alça na posição trancada""""'
$ws3.Range("A7").Value = 'Realocar'
$ws3.Range("A8").Value = 'Arrumar os 2 bebedouros que estão quebrados e vazando no P-2-S (remover o balde assim que concluído)'
$ws3.Range("A9").Value = 'Remover materias da doca do IB'
$ws3.Range("A10").Value = 'Fixar a placa do hidrante novamente na parede'
$ws3.Range("A11").Value = 'Fazer demarcações para totes vazios nas estações de single'
$ws3.Range("A12").Value = 'Fazer demarcação amarela para skutles de lixo'
$ws3.Range("A13").Value = 'Refoçar com os AAs a regra de estável, para que se o item cair, caia para dentro do bin. Ajustar posição P-2-S119M400'
$ws3.Range("A14").Value = 'Mover carrinhos para demarcação ou realizar demarcação dos carrinhos. Reforçar com AAs para manter os carrinhos na demarcação'
$ws3.Range("A15").Value = 'Fazer demarcações para as lixeiras da coluna FF05'
$ws3.Range("A16").Value = 'Demarcar buffer de carrinho de ICQA com fita de azul no P-2-S'
$ws3.Range("A17").Value = 'Remover os unifilas e cones próximo ao E-2 rua 139 altura do 200-285'
$ws3.Range("A18").Value = 'Destinar ventiladores'
$ws3.Range("A19").Value = 'Realocar escadas para local correto'
$ws3.Range("A20").Value = 'Paleteira movida para posição correta'
$ws3.Range("A21").Value = 'Verificar placas de piso molhado do Mod E e me enviar a planilha atualizada'
$ws3.Range("A22").Value = 'Cart removido e alocado na área correta.'
$ws3.Range("A23").Value = 'Colocar demarcações nas lixeiras dos armários e escritório ADM'
$ws3.Range("A24").Value = 'O elástico foi retirado'
$ws3.Range("A25").Value = 'Ensure cleaning routine or look for trash instalation'
$ws3.Range("A26").Value = 'Abrir ticket para a manutenção realizar o reparo'
$ws3.Range("A27").Value = 'Remover placas de """"Stow velocidade alta"""" no penultimo corredor proximo ao OB'
$ws3.Range("A28").Value = 'Define proper location to store it'
$ws3.Range("A29").Value = 'Fazer demarcação para carrinho de Rebin na estação.'
$ws3.Range("A30").Value = 'Orientar AAs do Pack a manter apenas 1 cart na estação.'
$ws3.Range("A31").Value = 'Mesa movida para area demarcada'
